# "Added the new button 13-17 on app"
# Populate the newly added buttons' Android->Rpi / Rpi columns (C/D) for
# rows 6-10 with the next sequence of test numbers (13-17), and update the
# testing notes/status for a handful of rows further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New button test numbers (13-17) for rows 6-10, columns C (Android -> Rpi)
# and D (Rpi) - these were previously blank.
$ws.Range("C6").Value = 13
$ws.Range("D6").Value = 13

$ws.Range("C7").Value = 14
$ws.Range("D7").Value = 14

$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 15

$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 16

$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 17

# Rows 18 & 19: status flips from the old error note to "OK", plus a new
# numeric reading of 6 in column H.
$ws.Range("G18").Value = "OK"
$ws.Range("H18").Value = 6

$ws.Range("G19").Value = "OK"
$ws.Range("H19").Value = 6

# Row 20: replace the old free-text note with a numeric reading of 12.
$ws.Range("H20").Value = 12

# Row 21: status note changed from "?" to the new phone-crash note.
$ws.Range("G21").Value = "PHONE CRASHES SOMETIMES "

# Row 11 note: LEDs are now reported as working fine.
$ws.Range("H11").Value = "Leds working Okay"

# Selection ends on D11, matching the author's last-touched cell.
$ws.Range("D11").Select()
